$wb = $excel.ActiveWorkbook
$students = $wb.Worksheets.Item("Students")

# Add "Errors" sheet right after "Students"
$errorsSheet = $wb.Worksheets.Add([System.Type]::Missing, $students)
$errorsSheet.Name = "Errors"

# Add "Warnings" sheet right after "Errors"
$warningsSheet = $wb.Worksheets.Add([System.Type]::Missing, $errorsSheet)
$warningsSheet.Name = "Warnings"

# Validation messages for the Errors sheet. Each literally begins with a
# single quote character, so we route the text through a formula
# ("=""...""") and then paste-special the computed value back over itself.
# This yields a plain shared-string text cell without Excel's "quote
# prefix" (leading-apostrophe) interpretation/styling kicking in.
$messages = @(
    "'Sheet ""Students"" Row: 2 Missing ""LAST NAME""',",
    "'Sheet ""Students"" Row: 2 Missing ""FIRST NAME""',",
    "'Sheet ""Students"" Row: 2 Missing ""STUDENT ID""',",
    "'Sheet ""Students"" Row: 2 Missing ""BIRTH DT""',",
    "'Sheet ""Students"" Row: 2 Missing ""OFF CLS""',",
    "'Sheet ""Students"" Row: 2 Invalid birthday """"',"
)

for ($i = 0; $i -lt $messages.Length; $i++) {
    $row = $i + 1
    $cell = $errorsSheet.Cells.Item($row, 1)
    $escaped = $messages[$i].Replace('"', '""')
    $cell.Formula = '="' + $escaped + '"'
    $cell.Copy()
    $cell.PasteSpecial(-4163)  # xlPasteValues
}
$excel.CutCopyMode = $false

$errorsSheet.Range("A1:A6").Select()

# Warnings sheet stays empty; make it the active/selected tab
$warningsSheet.Activate()
